# Recompute Psen1 -> Notch4 LR-pair stats per Dr Hou's advice (natmi rerun):
# the sender x target-cluster grid now spans all 5 clusters (ECs, FAPs, M1, M2, sCs)
# instead of only 4 target clusters, and ligand/receptor-expressing-cell counts rise
# from 1 to 3 (of the enlarged per-cluster cell pool), so every downstream NATMI score
# (detection rate, specificity, edge weights) is recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Psen1","Notch4","ECs",3,1,16.38931533333333,49.167946,0.1272611691689642,0.1272611691689643,3,1,47.253984,141.761952,0.8884015114449556,0.8884015114449555,774.460444532288,6970.144000790592,0.11305901503796,0.11305901503796),
    @("ECs","Psen1","Notch4","FAPs",3,1,16.38931533333333,49.167946,0.1272611691689642,0.1272611691689643,3,1,4.022517333333333,12.067552,0.07562559124637756,0.07562559124637755,65.92630500979911,593.336745088192,0.009624201161108197,0.009624201161108197),
    @("ECs","Psen1","Notch4","M1",3,1,16.38931533333333,49.167946,0.1272611691689642,0.1272611691689643,3,1,0.7893206666666668,2.367962,0.01483967305870774,0.01483967305870774,12.93642530511689,116.427827746052,0.001888514143536327,0.001888514143536327),
    @("ECs","Psen1","Notch4","M2",3,1,16.38931533333333,49.167946,0.1272611691689642,0.1272611691689643,2,0.6666666666666666,0.803095,2.409285,0.01509863828272948,0.01509863828272948,13.16217719762333,118.45959477861,0.001921470360719437,0.001921470360719437),
    @("ECs","Psen1","Notch4","sCs",3,1,16.38931533333333,49.167946,0.1272611691689642,0.1272611691689643,3,1,0.320979,0.9629369999999999,0.006034585967229564,0.006034585967229563,5.260626046377999,47.345634417402,0.0007679684656402593,0.0007679684656402593),
    @("FAPs","Psen1","Notch4","ECs",3,1,23.071008,69.213024,0.1791437526383466,0.1791437526383466,3,1,47.253984,141.761952,0.8884015114449556,0.8884015114449555,1090.197042895872,9811.773386062849,0.1591515806098284,0.1591515806098284),
    @("FAPs","Psen1","Notch4","FAPs",3,1,23.071008,69.213024,0.1791437526383466,0.1791437526383466,3,1,4.022517333333333,12.067552,0.07562559124637756,0.07562559124637755,92.803529577472,835.231766197248,0.01354785221136977,0.01354785221136977),
    @("FAPs","Psen1","Notch4","M1",3,1,23.071008,69.213024,0.1791437526383466,0.1791437526383466,3,1,0.7893206666666668,2.367962,0.01483967305870774,0.01483967305870774,18.210423415232,163.893810737088,0.002658434719663076,0.002658434719663076),
    @("FAPs","Psen1","Notch4","M2",3,1,23.071008,69.213024,0.1791437526383466,0.1791437526383466,2,0.6666666666666666,0.803095,2.409285,0.01509863828272948,0.01509863828272948,18.52821116976,166.75390052784,0.002704826721697161,0.002704826721697161),
    @("FAPs","Psen1","Notch4","sCs",3,1,23.071008,69.213024,0.1791437526383466,0.1791437526383466,3,1,0.320979,0.9629369999999999,0.006034585967229564,0.006034585967229563,7.405309076832,66.64778169148799,0.001081058375788211,0.00108105837578821),
    @("M1","Psen1","Notch4","ECs",3,1,37.292974,111.878922,0.2895757007844777,0.2895757007844777,3,1,47.253984,141.761952,0.8884015114449556,0.8884015114449555,1762.241596708416,15860.17437037575,0.2572594902546622,0.2572594902546622),
    @("M1","Psen1","Notch4","FAPs",3,1,37.292974,111.878922,0.2895757007844777,0.2895757007844777,3,1,4.022517333333333,12.067552,0.07562559124637756,0.07562559124637755,150.0116343265493,1350.104708938944,0.02189933358241025,0.02189933358241024),
    @("M1","Psen1","Notch4","M1",3,1,37.292974,111.878922,0.2895757007844777,0.2895757007844777,3,1,0.7893206666666668,2.367962,0.01483967305870774,0.01483967305870774,29.43611509966267,264.9250358969641,0.004297208725387828,0.004297208725387827),
    @("M1","Psen1","Notch4","M2",3,1,37.292974,111.878922,0.2895757007844777,0.2895757007844777,2,0.6666666666666666,0.803095,2.409285,0.01509863828272948,0.01509863828272948,29.94980095453,269.54820859077,0.004372198761612734,0.004372198761612733),
    @("M1","Psen1","Notch4","sCs",3,1,37.292974,111.878922,0.2895757007844777,0.2895757007844777,3,1,0.320979,0.9629369999999999,0.006034585967229564,0.006034585967229563,11.970261501546,107.732353513914,0.001747469460404676,0.001747469460404676),
    @("M2","Psen1","Notch4","ECs",3,1,43.83143633333333,131.494309,0.3403461170089362,0.3403461170089362,3,1,47.253984,141.761952,0.8884015114449556,0.8884015114449555,2071.209991192352,18640.88992073117,0.3023640047651606,0.3023640047651606),
    @("M2","Psen1","Notch4","FAPs",3,1,43.83143633333333,131.494309,0.3403461170089362,0.3403461170089362,3,1,4.022517333333333,12.067552,0.07562559124637756,0.07562559124637755,176.3127123957298,1586.814411561568,0.0257388763272096,0.02573887632720959),
    @("M2","Psen1","Notch4","M1",3,1,43.83143633333333,131.494309,0.3403461170089362,0.3403461170089362,3,1,0.7893206666666668,2.367962,0.01483967305870774,0.01483967305870774,34.59705854758423,311.373526928258,0.005050625103213304,0.005050625103213302),
    @("M2","Psen1","Notch4","M2",3,1,43.83143633333333,131.494309,0.3403461170089362,0.3403461170089362,2,0.6666666666666666,0.803095,2.409285,0.01509863828272948,0.01509863828272948,35.20080736211833,316.807266259065,0.005138762911649453,0.005138762911649452),
    @("M2","Psen1","Notch4","sCs",3,1,43.83143633333333,131.494309,0.3403461170089362,0.3403461170089362,3,1,0.320979,0.9629369999999999,0.006034585967229564,0.006034585967229563,14.068970602837,126.620735425533,0.002053847901703197,0.002053847901703197),
    @("sCs","Psen1","Notch4","ECs",3,1,8.200153666666667,24.600461,0.06367326039927532,0.06367326039927532,3,1,47.253984,141.761952,0.8884015114449556,0.8884015114449555,387.4899301622081,3487.409371459872,0.05656742077734443,0.05656742077734442),
    @("sCs","Psen1","Notch4","FAPs",3,1,8.200153666666667,24.600461,0.06367326039927532,0.06367326039927532,3,1,4.022517333333333,12.067552,0.07562559124637756,0.07562559124637755,32.98526026016356,296.867342341472,0.004815327964279754,0.004815327964279754),
    @("sCs","Psen1","Notch4","M1",3,1,8.200153666666667,24.600461,0.06367326039927532,0.06367326039927532,3,1,0.7893206666666668,2.367962,0.01483967305870774,0.01483967305870774,6.472550758942446,58.25295683048201,0.0009448903669072085,0.0009448903669072083),
    @("sCs","Psen1","Notch4","M2",3,1,8.200153666666667,24.600461,0.06367326039927532,0.06367326039927532,2,0.6666666666666666,0.803095,2.409285,0.01509863828272948,0.01509863828272948,6.585502408931667,59.269521680385,0.0009613795270507016,0.0009613795270507014),
    @("sCs","Psen1","Notch4","sCs",3,1,8.200153666666667,24.600461,0.06367326039927532,0.06367326039927532,3,1,0.320979,0.9629369999999999,0.006034585967229564,0.006034585967229563,2.632077123773,23.688694113957,0.0003842417636932208,0.0003842417636932207)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
